# Update Pertanggal 12 November 2020
# Fill in / correct DESCRIPTION (column C) values on the "List" sheet for a
# number of API-List rows, then leave the selection on the last-updated cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# transaction.delete.master.setBloodAglutinogenType -> capitalize description
$ws.Range("C13").Value = "Menghapus Data Jenis Golongan Darah"

# Previously-empty descriptions for the "initialize" API rows - now filled in
$ws.Range("C24").Value = "Menginisialisasi Data Hari Libur Kebijakan Pemerintah"
$ws.Range("C25").Value = "Menginisialisasi Data Hari Libur Nasional"
$ws.Range("C26").Value = "Menginisialisasi Data Model Barang"
$ws.Range("C27").Value = "Menginisialisasi Data Jenis Barang"
$ws.Range("C28").Value = "Menginisialisasi Data Periode"
$ws.Range("C29").Value = "Menginisialisasi Data Orang"
$ws.Range("C30").Value = "Menginisialisasi Data Akun EMail Orang"
$ws.Range("C31").Value = "Menginisialisasi Data Jenis Kelamin Orang"
$ws.Range("C32").Value = "Menginisialisasi Data Jenis Produk"
$ws.Range("C33").Value = "Menginisialisasi Data Agama"
$ws.Range("C34").Value = "Menginisialisasi Data Merk Dagang"

# transaction.read.master.getDataListPeriod had the wrong (copy/pasted) description
$ws.Range("C47").Value = "Mendapatkan Daftar Periode dari Data Master"

# transaction.undelete.master.setBloodAglutinogenType -> capitalize description
$ws.Range("C56").Value = "Membatalkan Penghapusan Data Jenis Golongan Darah"

# Move the active selection down to where the edits were made
$ws.Range("C57").Select()
